# BOM.xlsx edit script - "Docs & BOM: Excluded tax"
# Recomputes unit prices (Á price) to exclude tax, updates a distributor
# (DigiKey -> Farnell) for the LT3750 part, bumps the LEDEX solenoid
# quantity, tweaks column widths and the view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Unit price ("A price", column H) updates - tax excluded.
#    The Estimated cost column (I) holds formulas and recalculates
#    automatically.
# ---------------------------------------------------------------------
$ws.Range("H3").Value = 297
$ws.Range("H4").Value = 215.4
$ws.Range("H5").Value = 268.8
$ws.Range("H6").Value = 147.2
$ws.Range("H7").Value = 154.808
$ws.Range("H8").Value = 203.2
$ws.Range("H9").Value = 351.2
$ws.Range("H10").Value = 359.2
$ws.Range("H11").Value = 60.14
$ws.Range("H20").Value = 198.352
$ws.Range("H21").Value = 309
$ws.Range("H36").Value = 175.2

# ---------------------------------------------------------------------
# 2. LT3750 capacitor charging controller (row 11): switched distributor
#    from DigiKey to Farnell, with new distributor/manufacturer part
#    numbers.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "Farnell"
$ws.Range("E11").Value = "ANALOG DEVICES LT3750EMS#TRPBF"
$ws.Range("G11").Value = "LT3750EMS#TRPBF"

# ---------------------------------------------------------------------
# 3. LEDEX solenoid (row 34): quantity doubled 3 -> 6.
# ---------------------------------------------------------------------
$ws.Range("C34").Value = 6

# ---------------------------------------------------------------------
# 4. Hyperlinks: rebuild the hyperlink list so the D11 "display" text
#    matches the new distributor, keeping every other hyperlink intact.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D3"), "https://hobbyking.com/en_us/turnigy-multistar-4225-610kv-16pole-multi-rotor-outrunner.html", "", "", "HobbyKing ")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.rcflight.se/visaprodukt.aspx?id=2951&p=t-motor-mt2814-400-kv", "", "", "RCflight")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://uav-en.tmotor.com/", "", "", "Tiger motors")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.elefun.se/p/prod.aspx?v=54896", "", "", "elefun")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://hobbyking.com/en_us/aerostar-30a-rvs-g2-32bit-2-4s-electronic-speed-controller-w-reverse-function-4a-5-6v-sbec.html", "", "", "HobbyKing")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://hobbyking.com/en_us/turnigy-plush-32-30a-2-4s-brushless-speed-controller-w-bec-rev1-1-0.html", "", "", "HobbyKing")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.amazon.com/Rakstore-ESP32-DevKitC-VIE-ESP32-WROVER-IE-Development-Bluetooth/dp/B09BM2D6HJ", "", "", "Amazon")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.elefun.se/p/prod.aspx?v=63197", "", "", "Elefun")
$ws.Hyperlinks.Add($ws.Range("F9"), "http://www.gaonengbattery.com/", "", "", "GNB")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://droneit.se/product/tattu-r-line-version-4-0-1300mah-22-2v-130c-6s1p-lipo-battery-pack-with-xt60-plug/", "", "", "Droneit")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.digikey.se/sv/products/detail/analog-devices-inc/LT3750EMS-PBF/1619999", "", "", "Farnell")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.digikey.se/sv/supplier-centers/analog-devices", "", "", "Analog-devices")
$ws.Hyperlinks.Add($ws.Range("D19"), "https://www.symmetryelectronics.com/products/ic-haus/ic-px2604odfn8-3x3/", "", "", "Symmetry Electronics")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.ichaus.de/product/ic-px-series/", "documents", "", "iC-Haus")
$ws.Hyperlinks.Add($ws.Range("D20"), "https://www.mouser.se/ProductDetail/Same-Sky/AMT102-0512-I5000-S?qs=gTYE2QTfZfSxiIvKD%252BmReg%3D%3D", "", "", "Mouser")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.sameskydevices.com/", "", "", "Same sky")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://uk.farnell.com/broadcom-limited/aedb-9140-a13/encoder-3channel-500cpr-8mm/dp/1161087", "", "", "Farnell")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://uk.farnell.com/b/broadcom", "", "", "BROADCOM")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.we-online.com/en/components/products/WSEN-ISDS", "", "", "Würth Electronik")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://www.we-online.com/en", "", "", "Würth Elektronik")
$ws.Hyperlinks.Add($ws.Range("D23"), "https://www.mouser.se/ProductDetail/Texas-Instruments/OPT8241NBN?qs=cGEy3R83DS%2FxFMUAL%252BoBvw%3D%3D", "", "", "Mouser")
$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.mouser.se/manufacturer/texas-instruments/", "", "", "Texas Instruments")
$ws.Hyperlinks.Add($ws.Range("D34"), "https://uk.farnell.com/ledex/195207-228/solenoid-tubular-10w-25-91x52/dp/3996096", "", "", "Farnell")
$ws.Hyperlinks.Add($ws.Range("F34"), "https://uk.farnell.com/b/ledex", "", "", "LEDEX")
$ws.Hyperlinks.Add($ws.Range("D35"), "https://www.autodoc.se/as-pl/12111476", "", "", "Autodoc")
$ws.Hyperlinks.Add($ws.Range("F35"), "https://as-pl.com/en/main", "", "", "AS-PL")
$ws.Hyperlinks.Add($ws.Range("D36"), "https://www.elefun.se/p/prod.aspx?v=65193", "", "", "Elefun")
$ws.Hyperlinks.Add($ws.Range("F36"), "https://www.hobbywing.com/", "", "", "Hobbywing")

# ---------------------------------------------------------------------
# 5. Column widths tweaked slightly.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 71.42
$ws.Columns.Item(3).ColumnWidth = 14.37
$ws.Columns.Item(4).ColumnWidth = 18.96
$ws.Columns.Item(5).ColumnWidth = 34.41
$ws.Columns.Item(6).ColumnWidth = 15.63
$ws.Columns.Item(7).ColumnWidth = 115.8
$ws.Columns.Item(8).ColumnWidth = 7.14
$ws.Columns.Item(9).ColumnWidth = 13.54
$ws.Columns.Item(10).ColumnWidth = 34.96
$ws.Columns.Item(11).ColumnWidth = 40.53

# ---------------------------------------------------------------------
# 6. View state: zoom out to 75% and move the active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 75
$ws.Range("C40").Select()
